# Update "想去人数" (interested count) figures for two worksheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3,4,5 in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1096
$wsExpo.Range("F4").Value = 2516
$wsExpo.Range("F5").Value = 214

# Sheet "全部类型" (All types) - rows 5,6,8 in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1096
$wsAll.Range("F6").Value = 2516
$wsAll.Range("F8").Value = 214
